$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows (6-12), mirroring the existing "Bag" method rows ---
$dates = @(
    42601.767430555556,
    42601.76935185185,
    42601.770104166666,
    42601.771111111113,
    42601.771539351852,
    42601.772905092592,
    42601.7733912037
)
$cValues = @(23, 21, 21, 83, 26, 24, 24)

# Copy the date/time number format (and style) used by the existing A-column
# cells onto the new ones before writing values, so the new cells share
# style index 1 instead of minting a new custom number format.
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A6:A12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

for ($i = 0; $i -lt 7; $i++) {
    $r = 6 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = "Bag"
    $ws.Cells.Item($r, 3).Value = $cValues[$i]
    $ws.Cells.Item($r, 4).Value = 17
    $ws.Cells.Item($r, 5).Value = 2
    $ws.Cells.Item($r, 6).Value = 0
    $ws.Cells.Item($r, 7).Value = 1
    $ws.Cells.Item($r, 8).Value = 0
    $ws.Cells.Item($r, 9).Value = 100
    $ws.Cells.Item($r, 10).Value = 2
    $ws.Cells.Item($r, 11).Value = 0
    $ws.Cells.Item($r, 12).Value = 100
    $ws.Cells.Item($r, 13).Value = 0
}

# --- Widen column A slightly to fit the refreshed "best fit" data ---
$ws.Columns.Item(1).ColumnWidth = 14
